$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# "About" sheet: refresh the source citation (now MIT paper instead of
# Bloomberg New Energy Finance), drop the embedded chart image + the
# "as of" date stamp, and add a methodology note.
# ----------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Remove the date stamp in C1 entirely (value + its date-number style).
$wsAbout.Range("C1").Clear()

# New source info.
$wsAbout.Range("B3").Value = "Massachusetts Institute of Technology"
$wsAbout.Range("B4").Value = 2021
$wsAbout.Range("B5").Value = "Re-examining rates of lithium-ion battery technology improvement and cost decline"
$wsAbout.Range("B6").Value = "https://pubs.rsc.org/en/content/articlepdf/2021/ee/d0ee02681f?page=search"
$wsAbout.Range("B7").Value = "Abstract"

# Old footnote text removed, but keep the italic style on C8.
$wsAbout.Range("C8").ClearContents()

# New methodology note in row 9.
$wsAbout.Range("A9").Value = "Note: We take the average of learning rates quoted in the Abstract (20%-27%)"

# Remove the embedded picture (chart reproduction) from the sheet.
if ($wsAbout.Shapes.Count -gt 0) {
    for ($i = $wsAbout.Shapes.Count; $i -ge 1; $i--) {
        $wsAbout.Shapes.Item($i).Delete()
    }
}

# ----------------------------------------------------------------------
# "PDiBCpDoC" sheet: the decline rate is now computed (average of the
# 20%-27% range from the MIT paper's abstract) instead of a hardcoded
# literal.
# ----------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("PDiBCpDoC")
$wsData.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"
$wsData.Range("A2").Value = "Batteries"
$wsData.Range("B2").Formula = "=AVERAGE(0.2,0.27)"

# ----------------------------------------------------------------------
# Selection / active-cell bookkeeping (matches the saved workbook state).
# Activate "PDiBCpDoC" first so its selection is recorded, then finish on
# "About" so that sheet remains the active tab, as in the original file.
# ----------------------------------------------------------------------
[void]$wsData.Activate()
[void]$wsData.Range("I4").Select()

[void]$wsAbout.Activate()
[void]$wsAbout.Range("A10").Select()
